$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.221.17'
$ws.Range("E2").Value = '  -1.13%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.574.11'
$ws.Range("E3").Value = '  -0.40%  '

# Row 4
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.73'
$ws.Range("E5").Value = '  +0.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.489'
$ws.Range("E6").Value = '  -1.87%  '

# Row 7
$ws.Range("E7").Value = '  -0.20%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.21'
$ws.Range("E8").Value = '  +0.24%  '

# Row 9
$ws.Range("E9").Value = '  -0.74%  '

# Row 10
$ws.Range("E10").Value = '  +0.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0866'
$ws.Range("E11").Value = '  +0.01%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.798.09'
$ws.Range("E12").Value = '  -0.53%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.577.24'
$ws.Range("E13").Value = '  -0.23%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.78'
$ws.Range("E14").Value = '  -1.13%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").Value = '  -0.81%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.219.20'
$ws.Range("E16").Value = '  -1.22%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.10'
$ws.Range("E17").Value = '  -1.28%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.36'
$ws.Range("E18").Value = '  +1.15%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.33'
$ws.Range("E19").Value = '  -0.30%  '

# Row 20
$ws.Range("E20").Value = '  -0.52%  '

# Row 21
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +0.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.44'
$ws.Range("E23").Value = '  -2.43%  '

# Row 24
$ws.Range("E24").Value = '  +0.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.01'
$ws.Range("E25").Value = '  -0.47%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.67'
$ws.Range("E26").Value = '  -3.70%  '

# Row 27
$ws.Range("E27").Value = '  -0.44%  '

# Row 28
$ws.Range("E28").Value = '  -0.16%  '

# Row 29
$ws.Range("E29").Value = '  -1.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.13'
$ws.Range("E30").Value = '  -2.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0464'
$ws.Range("E31").Value = '  -1.49%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.18'
$ws.Range("E32").Value = '  -1.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.407.84'
$ws.Range("E33").Value = '  +2.74%  '

# Row 34
$ws.Range("E34").Value = '  -0.87%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +2.27%  '

# Row 36
$ws.Range("E36").Value = '  -1.01%  '

# Row 37
$ws.Range("E37").Value = '  -2.80%  '

# Row 38
$ws.Range("E38").Value = '  -1.78%  '

# Row 39
$ws.Range("E39").Value = '  -0.13%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.517'
$ws.Range("E40").Value = '  -2.58%  '

# Row 41
$ws.Range("E41").Value = '  -0.11%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.83'
$ws.Range("E43").Value = '  +4.43%  '

# Row 44
$ws.Range("E44").Value = '  +2.30%  '

# Row 45
$ws.Range("E45").Value = '  +0.79%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.83'
$ws.Range("E46").Value = '  -0.29%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.710.39'
$ws.Range("E47").Value = '  -0.65%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.16'
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₇0993'
$ws.Range("E49").Value = '  -0.79%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("E50").Value = '  +0.08%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0494'
$ws.Range("E51").Value = '  -0.02%  '
